# Update "想去人数" (column F) figures across the four worksheets to match
# the newly regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 1290
$ws.Range("F7").Value = 377
$ws.Range("F8").Value = 1304
$ws.Range("F9").Value = 909
$ws.Range("F11").Value = 200
$ws.Range("F16").Value = 3020
$ws.Range("F17").Value = 2659
$ws.Range("F19").Value = 31
$ws.Range("F21").Value = 325
$ws.Range("F22").Value = 242
$ws.Range("F24").Value = 5409
$ws.Range("F25").Value = 595
$ws.Range("F30").Value = 1138
$ws.Range("F31").Value = 73
$ws.Range("F32").Value = 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 236
$ws.Range("F9").Value = 40
$ws.Range("F23").Value = 326
$ws.Range("F25").Value = 3997
$ws.Range("F33").Value = 38

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1075
$ws.Range("F9").Value = 1368

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 1075
$ws.Range("F8").Value = 1368
$ws.Range("F13").Value = 1290
$ws.Range("F14").Value = 377
$ws.Range("F15").Value = 909
$ws.Range("F18").Value = 200
$ws.Range("F22").Value = 3020
$ws.Range("F23").Value = 2659
$ws.Range("F24").Value = 31
$ws.Range("F25").Value = 325
$ws.Range("F26").Value = 40
$ws.Range("F27").Value = 242
$ws.Range("F29").Value = 5409
$ws.Range("F30").Value = 595
$ws.Range("F40").Value = 326
$ws.Range("F42").Value = 1138
$ws.Range("F43").Value = 73
$ws.Range("F47").Value = 75
